$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 293, shifting existing rows 293:338 down to 294:339
$ws.Rows.Item(293).Insert()

# Populate the newly inserted row 293 with the new weekly record (same as the
# following week's row, but with updated date/volume/price figures).
$ws.Cells.Item(293, 1).Value = 11
$ws.Cells.Item(293, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(293, 3).Value = "Bíobío"
$ws.Cells.Item(293, 4).Value = 45218
$ws.Cells.Item(293, 5).Value = 8
$ws.Cells.Item(293, 6).Value = 100112003
$ws.Cells.Item(293, 7).Value = "Ajo"
$ws.Cells.Item(293, 8).Value = "Chino"
$ws.Cells.Item(293, 9).Value = "Primera"
$ws.Cells.Item(293, 10).Value = 200
$ws.Cells.Item(293, 11).Value = 20000
$ws.Cells.Item(293, 12).Value = 21000
$ws.Cells.Item(293, 13).Value = 20500
$ws.Cells.Item(293, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(293, 15).Value = "China"
$ws.Cells.Item(293, 16).Value = 2050
$ws.Cells.Item(293, 17).Value = 10
$ws.Cells.Item(293, 18).Value = "Hortaliza"
